$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-10-23 Thursday"; New = "2025-10-24 Friday" },
    @{ Old = "704×5="; New = "673×8=" },
    @{ Old = "220×3="; New = "886×3=" },
    @{ Old = "886×5="; New = "760×6=" },
    @{ Old = "110×8="; New = "953×7=" },
    @{ Old = "687×3="; New = "962×2=" },
    @{ Old = "942×2="; New = "450×4=" },
    @{ Old = "638×7="; New = "562×5=" },
    @{ Old = "123×6="; New = "894×9=" },
    @{ Old = "809×7="; New = "938×5=" },
    @{ Old = "901×6="; New = "662×3=" },
    @{ Old = "687×3="; New = "223×9=" },
    @{ Old = "991×9="; New = "589×8=" },
    @{ Old = "484×2="; New = "941×9=" },
    @{ Old = "624×4="; New = "908×3=" },
    @{ Old = "133×4="; New = "456×4=" },
    @{ Old = "216×8="; New = "353×5=" },
    @{ Old = "278×5="; New = "853×3=" },
    @{ Old = "246×6="; New = "170×2=" },
    @{ Old = "480×2="; New = "186×6=" },
    @{ Old = "731×8="; New = "306×2=" },
    @{ Old = "238×2="; New = "412×7=" },
    @{ Old = "522×6="; New = "313×4=" },
    @{ Old = "559×5="; New = "941×5=" },
    @{ Old = "963×2="; New = "434×6=" },
    @{ Old = "686×5="; New = "360×8=" }
)

$searchStart = 0

foreach ($rep in $replacements) {
    $docEnd = $d.Content.End
    $r = $d.Range($searchStart, $docEnd)
    $found = $r.Find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.New, 1)
    $searchStart = $r.End
}
